$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("GET_Tests")
$ws2 = $wb.Worksheets.Item("POST Tests")

# Update the timeout values in column L (seconds -> milliseconds: 2 -> 750)
$ws1.Range("L1").Value = 750
$ws1.Range("L2").Value = 750
$ws1.Range("L3").Value = 750

$ws2.Range("L1").Value = 750
$ws2.Range("L2").Value = 750
$ws2.Range("L3").Value = 750

# Row heights on the POST Tests sheet become explicit/custom
$ws2.Rows.Item(2).RowHeight = 22.1
$ws2.Rows.Item(3).RowHeight = 24.35

# Selection on GET_Tests moves to L3 before switching away from it
$ws1.Activate()
$ws1.Range("L3").Select()

# POST Tests becomes the active sheet, scrolled/selected back at A1
$ws2.Activate()
$ws2.Range("A1").Select()
